$d = $word.ActiveDocument

# Locate the target paragraph: the Spanish "Detectar mediante..." quote paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like '*Detectar mediante*') {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$pPrCommon = '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:space="0" w:sz="0" w:val="nil"/><w:left w:space="0" w:sz="0" w:val="nil"/><w:bottom w:space="0" w:sz="0" w:val="nil"/><w:right w:space="0" w:sz="0" w:val="nil"/><w:between w:space="0" w:sz="0" w:val="nil"/></w:pBdr><w:shd w:fill="auto" w:val="clear"/><w:ind w:left="600" w:hanging="360"/><w:rPr><w:u w:val="none"/></w:rPr>'

$spanish = 'Detectar mediante lo que yo llamo "Gestures" el Propósito de una Interacción de un usuario (persona o servicio: flujos de navegación) en un Contexto dado, según los Datos del "diálogo" en un protocolo que permita inferir y facilitar la intención o el objetivo de la Interacción mediante "Suggestions".'
$english = 'Detect by means of user "Gestures" (person or service browsing flows) the Purpose of an Interaction, in a given Context, following Data of a of a "dialog" in a protocol such that the Interaction intention or objectives may be inferred and guided by means of "Suggestions".'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    '<w:p><w:pPr>' + $pPrCommon + '</w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + $spanish + '</w:t></w:r></w:p>' +
    '<w:p><w:pPr>' + $pPrCommon + '</w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + $english + '</w:t></w:r></w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($xml)

Write-Output ("Paragraph count now: " + $d.Paragraphs.Count)
